# Add two new rows (setPath / clearPath) to the "addt'l-functions" sheet
# and make that sheet the active tab (it was "phylip-programs" before).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("addt'l-functions")

# Clone the formatting of the last existing data row (row 15) onto the two
# new rows so the new cells pick up the same styles (left-aligned text,
# date number format for column B, etc.) instead of the engine's defaults.
$ws.Range("A15:D15").Copy()
$ws.Range("A16:D16").PasteSpecial(-4122)
$ws.Range("A15:D15").Copy()
$ws.Range("A17:D17").PasteSpecial(-4122)

# Row 16: setPath
$ws.Range("A16").Value = "setPath"
$ws.Range("B16").Value = 41629
$ws.Range("C16").Value = "No"
$ws.Range("D16").Value = "Sets path to the folder containing PHYLIP executables for current R session."

# Row 17: clearPath
$ws.Range("A17").Value = "clearPath"
$ws.Range("B17").Value = 41629
$ws.Range("C17").Value = "No"
$ws.Range("D17").Value = "Clears path to PHYLIP executables."

# Make "addt'l-functions" the active/selected sheet (was "phylip-programs").
$ws.Activate()
